# Update the apple_finrep output sheet: files/periods moved from 0222 (Feb)
# to 0322 (Mar), currency sort order rearranged (RO,PE,HU,EU,MX,LL,BG,BR,CA,
# CZ,CL,CO,NZ,AU,CH,NO,US,DK,PL,SE,JP,GB) and r_count/sum/built_in_total
# values refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  File = "87811004_0322_RO"; RCount = 423;  Currency = "RON"; Sum = 9983.42 }
    @{ Row = 3;  File = "87811004_0322_PE"; RCount = 15;   Currency = "PEN"; Sum = 141.96 }
    @{ Row = 4;  File = "87811004_0322_HU"; RCount = 954;  Currency = "HUF"; Sum = 1882603 }
    @{ Row = 5;  File = "87811004_0322_EU"; RCount = 716;  Currency = "EUR"; Sum = 3090.3 }
    @{ Row = 6;  File = "87811004_0322_MX"; RCount = 75;   Currency = "MXN"; Sum = 4825.1 }
    @{ Row = 7;  File = "87811004_0322_LL"; RCount = 37;   Currency = "USD"; Sum = 96.6 }
    @{ Row = 8;  File = "87811004_0322_BG"; RCount = 6;    Currency = "BGN"; Sum = 22.44 }
    @{ Row = 9;  File = "87811004_0322_BR"; RCount = 37;   Currency = "BRL"; Sum = 334.81 }
    @{ Row = 10; File = "87811004_0322_CA"; RCount = 324;  Currency = "CAD"; Sum = 1753.5 }
    @{ Row = 11; File = "87811004_0322_CZ"; RCount = 12;   Currency = "CZK"; Sum = 636.36 }
    @{ Row = 12; File = "87811004_0322_CL"; RCount = 31;   Currency = "CLP"; Sum = 63087 }
    @{ Row = 13; File = "87811004_0322_CO"; RCount = 31;   Currency = "COP"; Sum = 286230 }
    @{ Row = 14; File = "87811004_0322_NZ"; RCount = 51;   Currency = "NZD"; Sum = 208.33 }
    @{ Row = 15; File = "87811004_0322_AU"; RCount = 367;  Currency = "AUD"; Sum = 3178.64 }
    @{ Row = 16; File = "87811004_0322_CH"; RCount = 79;   Currency = "CHF"; Sum = 243.39 }
    @{ Row = 17; File = "87811004_0322_NO"; RCount = 26;   Currency = "NOK"; Sum = 682.5 }
    @{ Row = 18; File = "87811004_0322_US"; RCount = 1487; Currency = "USD"; Sum = 10299.1 }
    @{ Row = 19; File = "87811004_0322_DK"; RCount = 21;   Currency = "DKK"; Sum = 268.24 }
    @{ Row = 20; File = "87811004_0322_PL"; RCount = 46;   Currency = "PLN"; Sum = 517.94 }
    @{ Row = 21; File = "87811004_0322_SE"; RCount = 32;   Currency = "SEK"; Sum = 931.77 }
    @{ Row = 22; File = "87811004_0322_JP"; RCount = 23;   Currency = "JPY"; Sum = 5950 }
    @{ Row = 23; File = "87811004_0322_GB"; RCount = 477;  Currency = "GBP"; Sum = 1593.33 }
)

# Column E ("built_in_total") stores the same number as text in the source
# data. Force text storage up front (Excel would otherwise coerce a
# numeric-looking string back into a number), write the values, then
# restore the default "Normal" style so no stray number format sticks
# around on the cells.
$ws.Range("E2:E23").NumberFormat = "@"

foreach ($item in $data) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.File
    $ws.Cells.Item($r, 2).Value = $item.RCount
    $ws.Cells.Item($r, 3).Value = $item.Currency
    $ws.Cells.Item($r, 4).Value = $item.Sum
    $ws.Cells.Item($r, 5).Value = [string]$item.Sum
}

$ws.Range("E2:E23").Style = "Normal"
